$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the "methods" row (row 5) comma-separated value into separate cells
$ws.Range("B5").Value = "Topic modeling"
$ws.Range("C5").Value = "sentiment analysis"
$ws.Range("D5").Value = "mixed-methods"

# Split the "themes" row (row 6) comma-separated value into separate cells
$ws.Range("B6").Value = "Sosiale medier"
$ws.Range("C6").Value = "tekstanalyse"
$ws.Range("D6").Value = "misinformasjon"

# Touch column B width so it becomes an explicit custom width (matches diff)
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Update selection to the last edited cell
$ws.Range("D5").Select()
